# Insert a new weekly price observation as row 47 on the single data sheet,
# pushing the existing rows 47:116 down to 48:117 (dimension grows to A1:R117).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("47:47").Insert()

$ws.Range("A47").Value2 = 10
$ws.Range("B47").Value2 = "Vega Modelo de Temuco"
$ws.Range("C47").Value2 = "La Araucanía"
$ws.Range("D47").Value2 = 44792
$ws.Range("E47").Value2 = 9
$ws.Range("F47").Value2 = 100112035
$ws.Range("G47").Value2 = "Bruselas (repollito)"
$ws.Range("H47").Value2 = "Sin especificar"
$ws.Range("I47").Value2 = "Primera"
$ws.Range("J47").Value2 = 35
$ws.Range("K47").Value2 = 25000
$ws.Range("L47").Value2 = 25000
$ws.Range("M47").Value2 = 25000
$ws.Range("N47").Value2 = "$/malla 10 kilos"
$ws.Range("O47").Value2 = "Provincia de Quillota"
$ws.Range("P47").Value2 = 2500
$ws.Range("Q47").Value2 = 10
$ws.Range("R47").Value2 = "Hortaliza"
